$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.011.99"
$ws.Range("E2").Value = "  -0.94%  "

# Row 3
$ws.Range("D3").Value = "'3.053.59"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'583.00"
$ws.Range("E5").Value = "  -1.50%  "

# Row 6
$ws.Range("D6").Value = "'151.46"
$ws.Range("E6").Value = "  -2.44%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  -1.72%  "

# Row 9
$ws.Range("D9").Value = "'3.054.80"
$ws.Range("E9").Value = "  -0.97%  "

# Row 10
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  -2.96%  "

# Row 11
$ws.Range("D11").Value = "'5.84"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -2.30%  "

# Row 13
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -3.09%  "

# Row 14
$ws.Range("D14").Value = "'36.15"
$ws.Range("E14").Value = "  -3.86%  "

# Row 15
$ws.Range("E15").Value = "  +1.89%  "

# Row 16
$ws.Range("D16").Value = "'3.559.04"
$ws.Range("E16").Value = "  -1.22%  "

# Row 17
$ws.Range("D17").Value = "'7.14"
$ws.Range("E17").Value = "  -0.90%  "

# Row 18
$ws.Range("D18").Value = "'63.043.31"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19
$ws.Range("D19").Value = "'3.055.98"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("D20").Value = "'480.39"
$ws.Range("E20").Value = "  +0.86%  "

# Row 21
$ws.Range("D21").Value = "'14.31"
$ws.Range("E21").Value = "  -2.66%  "

# Row 22
$ws.Range("D22").Value = "'0.708"
$ws.Range("E22").Value = "  -1.79%  "

# Row 23
$ws.Range("D23").Value = "'7.52"
$ws.Range("E23").Value = "  -0.88%  "

# Row 24
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  -0.70%  "

# Row 25
$ws.Range("D25").Value = "'81.93"
$ws.Range("E25").Value = "  +0.78%  "

# Row 26
$ws.Range("D26").Value = "'12.67"
$ws.Range("E26").Value = "  -2.10%  "

# Row 27
$ws.Range("D27").Value = "'10.51"
$ws.Range("E27").Value = "  +4.94%  "

# Row 28
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("D29").Value = "'7.39"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("E31").Value = "  -1.46%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  +0.49%  "

# Row 33
$ws.Range("D33").Value = "'27.80"
$ws.Range("E33").Value = "  +1.75%  "

# Row 34
$ws.Range("E34").Value = "  -2.66%  "

# Row 35
$ws.Range("E35").Value = "  +0.82%  "

# Row 36
$ws.Range("D36").Value = "'0.0₃0813"
$ws.Range("E36").Value = "  -4.52%  "

# Row 37
$ws.Range("D37").Value = "'5.92"
$ws.Range("E37").Value = "  -3.43%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.20"
$ws.Range("E38").Value = "  -1.28%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  -6.34%  "

# Row 40
$ws.Range("D40").Value = "'9.20"
$ws.Range("E40").Value = "  -1.66%  "

# Row 41
$ws.Range("D41").Value = "'50.39"
$ws.Range("E41").Value = "  -0.89%  "

# Row 42
$ws.Range("D42").Value = "'428.21"
$ws.Range("E42").Value = "  -3.67%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.116"
$ws.Range("E43").Value = "  +3.82%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.287"
$ws.Range("E44").Value = "  +0.44%  "

# Row 45
$ws.Range("D45").Value = "'0.0362"
$ws.Range("E45").Value = "  -0.54%  "

# Row 46
$ws.Range("D46").Value = "'2.846.00"
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("D47").Value = "'38.02"
$ws.Range("E47").Value = "  -5.34%  "

# Row 48
$ws.Range("D48").Value = "'127.05"

# Row 49
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("D50").Value = "'25.20"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51
$ws.Range("E51").Value = "  -1.17%  "
